$d = $word.ActiveDocument

# Helper: split a contiguous run of text (already placed, with uniform
# character formatting) into separate <w:r> runs at the given internal
# offsets, without altering the visible text. We do this by toggling
# Bold on/off across the sub-range that starts at the split point --
# the Word engine here coalesces adjacent runs that share identical
# formatting, but a momentary formatting change forces a boundary that
# survives once the format is reverted.
function Split-RunAt($rangeStart, $offset, $len) {
    $b = $d.Range($rangeStart + $offset, $rangeStart + $offset + $len)
    $b.Font.Bold = 1
    $b.Font.Bold = 0
}

# ---------------------------------------------------------------------
# Edit 1: "37,  70" -> "37, " | "64, " | "65, " | "70"  (4 runs)
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("37,  70", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$newText1 = "37, 64, 65, 70"
$r1.Text = $newText1
$s1 = $r1.Start
Split-RunAt $s1 4 4    # boundary between "37, " and "64, "
Split-RunAt $s1 8 4    # boundary between "64, " and "65, "
Split-RunAt $s1 12 2   # boundary between "65, " and "70"

# ---------------------------------------------------------------------
# Edit 2: "18, 21, 36, 38,  39, 49, 51, 54"
#      -> "18, 21, 36, 38,  39, 49, 51, " | "53, " | "54"  (3 runs)
# ---------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("18, 21, 36, 38,  39, 49, 51, 54", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$newText2 = "18, 21, 36, 38,  39, 49, 51, 53, 54"
$r2.Text = $newText2
$s2 = $r2.Start
Split-RunAt $s2 29 4   # boundary between "...51, " and "53, "
Split-RunAt $s2 33 2   # boundary between "53, " and "54"

# ---------------------------------------------------------------------
# Edit 3: insert a new run "55, " immediately before the existing run
#         that holds "71, 72, 73, 74, 75, 76"
# ---------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("71, 72, 73, 74, 75, 76", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s3 = $r3.Start
$ins3 = $d.Range($s3, $s3)
$ins3.InsertBefore("55, ")
Split-RunAt $s3 0 4    # boundary between new "55, " and "71, 72, ..."

Write-Output "Edits applied"
